$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 645.7
$ws.Range("I2").Value = 645.7
$ws.Range("K2").Value = 645.7
$ws.Range("M2").Value = -532.7
$ws.Range("H11").Value = 188.28572
$ws.Range("I11").Value = 188.28572
$ws.Range("K11").Value = 188.28572
$ws.Range("M11").Value = -48.28572
$ws.Range("H40").Value = 1846.6086
$ws.Range("I40").Value = 1870.5405
$ws.Range("J40").Value = 1748.2222
$ws.Range("K40").Value = 1870.5405
$ws.Range("L40").Value = 1748.2222
$ws.Range("M40").Value = -1695.5405
$ws.Range("N40").Value = -2098.2222
$ws.Range("H51").Value = 2333.3333
$ws.Range("J51").Value = 3200
$ws.Range("L51").Value = 3200
$ws.Range("N51").Value = -4168
$ws.Range("H53").Value = 44095.652
$ws.Range("I53").Value = 91930.55
$ws.Range("J53").Value = 247
$ws.Range("K53").Value = 91930.55
$ws.Range("L53").Value = 247
$ws.Range("M53").Value = -91293.55
$ws.Range("N53").Value = -1521
$ws.Range("H106").Value = 121217350
$ws.Range("J106").Value = 333338340
$ws.Range("L106").Value = 333338340
$ws.Range("N106").Value = -333339602
$ws.Range("H133").Value = 24375
$ws.Range("J133").Value = 24375
$ws.Range("L133").Value = 24375
$ws.Range("N133").Value = -34495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2123.7778
$ws.Range("I94").Value = 1729.8182
$ws.Range("J94").Value = 2742.8572
$ws.Range("K94").Value = 1729.8182
$ws.Range("L94").Value = 2742.8572
$ws.Range("M94").Value = -1278.8182
$ws.Range("N94").Value = -3644.8572
$ws.Range("H107").Value = 1284.4615
$ws.Range("I107").Value = 1349.875
$ws.Range("J107").Value = 1179.8
$ws.Range("K107").Value = 1349.875
$ws.Range("L107").Value = 1179.8
$ws.Range("M107").Value = 570.125
$ws.Range("N107").Value = -5019.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2367.744
$ws.Range("I31").Value = 1446.8182
$ws.Range("J31").Value = 2705.4167
$ws.Range("K31").Value = 1446.8182
$ws.Range("L31").Value = 2705.4167
$ws.Range("M31").Value = -1151.8182
$ws.Range("N31").Value = -3295.4167
$ws.Range("H34").Value = 2367.744
$ws.Range("I34").Value = 1446.8182
$ws.Range("J34").Value = 2705.4167
$ws.Range("K34").Value = 1446.8182
$ws.Range("L34").Value = 2705.4167
$ws.Range("M34").Value = -1244.8182
$ws.Range("N34").Value = -3109.4167
$ws.Range("H107").Value = 515.90247
$ws.Range("I107").Value = 416.08334
$ws.Range("J107").Value = 656.82355
$ws.Range("K107").Value = 416.08334
$ws.Range("L107").Value = 656.82355
$ws.Range("M107").Value = 1503.91666
$ws.Range("N107").Value = -4496.82355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2662.0886
$ws.Range("I68").Value = 3898.0303
$ws.Range("J68").Value = 1775.4348
$ws.Range("K68").Value = 11694.0909
$ws.Range("L68").Value = 5326.3044
$ws.Range("M68").Value = -10883.0909
$ws.Range("N68").Value = -6948.3044
$ws.Range("H71").Value = 2662.0886
$ws.Range("I71").Value = 3898.0303
$ws.Range("J71").Value = 1775.4348
$ws.Range("K71").Value = 35082.2727
$ws.Range("L71").Value = 15978.9132
$ws.Range("M71").Value = -31026.2727
$ws.Range("N71").Value = -24090.9132
$ws.Range("H107").Value = 1101.2106
$ws.Range("I107").Value = 339.4
$ws.Range("J107").Value = 1263.2979
$ws.Range("K107").Value = 1018.2
$ws.Range("L107").Value = 3789.8937
$ws.Range("M107").Value = 901.8000000000001
$ws.Range("N107").Value = -7629.893700000001
$ws.Range("H122").Value = 464.85715
$ws.Range("I122").Value = 453.16666
$ws.Range("J122").Value = 535
$ws.Range("K122").Value = 4078.49994
$ws.Range("L122").Value = 4815
$ws.Range("M122").Value = -1628.49994
$ws.Range("N122").Value = -9715
$ws.Range("H125").Value = 3556.8462
$ws.Range("J125").Value = 3753.25
$ws.Range("L125").Value = 11259.75
$ws.Range("N125").Value = -21099.75
$ws.Range("H131").Value = 12360728
$ws.Range("I131").Value = 5882845.5
$ws.Range("J131").Value = 13890228
$ws.Range("K131").Value = 17648536.5
$ws.Range("L131").Value = 41670684
$ws.Range("M131").Value = -17643496.5
$ws.Range("N131").Value = -41680764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 999.3333
$ws.Range("I97").Value = 999.3333
$ws.Range("K97").Value = 999.3333
$ws.Range("M97").Value = -503.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 125003750
$ws.Range("I40").Value = 166670080
$ws.Range("J40").Value = 4749.5
$ws.Range("K40").Value = 166670080
$ws.Range("L40").Value = 4749.5
$ws.Range("M40").Value = -166669944
$ws.Range("N40").Value = -5021.5
$ws.Range("H46").Value = 25642114
$ws.Range("I46").Value = 33334294
$ws.Range("J46").Value = 1516.6666
$ws.Range("K46").Value = 33334294
$ws.Range("L46").Value = 1516.6666
$ws.Range("M46").Value = -33334106
$ws.Range("N46").Value = -1892.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H81").Value = 1671.1666
$ws.Range("I81").Value = 1666.6666
$ws.Range("J81").Value = 1675.6666
$ws.Range("K81").Value = 3333.3332
$ws.Range("L81").Value = 3351.3332
$ws.Range("M81").Value = -2272.3332
$ws.Range("N81").Value = -5473.3332
$ws.Range("H84").Value = 1671.1666
$ws.Range("I84").Value = 1666.6666
$ws.Range("J84").Value = 1675.6666
$ws.Range("K84").Value = 16666.666
$ws.Range("L84").Value = 16756.666
$ws.Range("M84").Value = -11362.666
$ws.Range("N84").Value = -27364.666
$ws.Range("H96").Value = 2503.5
$ws.Range("J96").Value = 3056
$ws.Range("L96").Value = 3056
$ws.Range("N96").Value = -5802
$ws.Range("H113").Value = 940
$ws.Range("I113").Value = 695.12
$ws.Range("K113").Value = 2085.36
$ws.Range("M113").Value = 84.63999999999987
$ws.Range("H126").Value = 1057.4762
$ws.Range("I126").Value = 791.5833
$ws.Range("J126").Value = 1412
$ws.Range("K126").Value = 2374.7499
$ws.Range("L126").Value = 4236
$ws.Range("M126").Value = 95.2501000000002
$ws.Range("N126").Value = -9176
$ws.Range("H136").Value = 974.03125
$ws.Range("I136").Value = 710.64
$ws.Range("J136").Value = 1914.7142
$ws.Range("K136").Value = 2131.92
$ws.Range("L136").Value = 5744.142599999999
$ws.Range("M136").Value = 418.0799999999999
$ws.Range("N136").Value = -10844.1426
